# jason and wooden plane have more stable aoa
# Update the "Wooden Plane" AoA (angle of attack) / Lift Coefficient table
# on Sheet1 (column B, rows 64-88) to the new, more stable lift-coefficient
# curve values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B64").Value = -0.72
$ws.Range("B65").Value = -0.92
$ws.Range("B66").Value = -1.02
$ws.Range("B67").Value = -1.04
$ws.Range("B68").Value = -1.02
$ws.Range("B69").Value = -0.98
$ws.Range("B70").Value = -0.94
$ws.Range("B71").Value = -0.82
$ws.Range("B72").Value = -0.68
$ws.Range("B73").Value = -0.57
$ws.Range("B74").Value = -0.4
$ws.Range("B75").Value = -0.2
$ws.Range("B77").Value = 0.2
$ws.Range("B78").Value = 0.4
$ws.Range("B79").Value = 0.57
$ws.Range("B80").Value = 0.68
$ws.Range("B81").Value = 0.82
$ws.Range("B82").Value = 0.94
$ws.Range("B83").Value = 0.98
$ws.Range("B84").Value = 1.02
$ws.Range("B85").Value = 1.04
$ws.Range("B86").Value = 1.02
$ws.Range("B87").Value = 0.92
$ws.Range("B88").Value = 0.71

# Restore the view/selection state that was active when the author saved:
# scrolled to the top of the "Wooden Plane" table (column A) with B65 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B65").Select()
